{"js": "// Replace each two-digit multiplication expression in the document's\n// table cells with its new value, as described by the diff.\nconst replacements = [\n  [\"21\u00d736=\", \"77\u00d779=\"],\n  [\"72\u00d789=\", \"94\u00d793=\"],\n  [\"82\u00d788=\", \"86\u00d731=\"],\n  [\"36\u00d773=\", \"91\u00d723=\"],\n  [\"54\u00d763=\", \"31\u00d772=\"],\n  [\"59\u00d724=\", \"69\u00d723=\"],\n  [\"17\u00d718=\", \"33\u00d750=\"],\n  [\"91\u00d748=\", \"64\u00d776=\"],\n  [\"32\u00d739=\", \"82\u00d720=\"],\n  [\"78\u00d765=\", \"31\u00d715=\"],\n  [\"95\u00d747=\", \"73\u00d745=\"],\n  [\"36\u00d718=\", \"58\u00d770=\"],\n  [\"26\u00d732=\", \"23\u00d739=\"],\n  [\"78\u00d754=\", \"86\u00d748=\"],\n  [\"58\u00d726=\", \"56\u00d750=\"],\n  [\"15\u00d762=\", \"57\u00d711=\"],\n  [\"79\u00d740=\", \"65\u00d763=\"],\n  [\"53\u00d778=\", \"24\u00d797=\"],\n  [\"98\u00d729=\", \"31\u00d738=\"],\n  [\"20\u00d795=\", \"29\u00d770=\"],\n  [\"74\u00d750=\", \"27\u00d797=\"],\n  [\"44\u00d731=\", \"60\u00d788=\"],\n  [\"26\u00d753=\", \"81\u00d763=\"],\n  [\"54\u00d748=\", \"58\u00d767=\"],\n  [\"59\u00d799=\", \"12\u00d789=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"21\u00d736=\", \"77\u00d779=\"),\n  @(\"72\u00d789=\", \"94\u00d793=\"),\n  @(\"82\u00d788=\", \"86\u00d731=\"),\n  @(\"36\u00d773=\", \"91\u00d723=\"),\n  @(\"54\u00d763=\", \"31\u00d772=\"),\n  @(\"59\u00d724=\", \"69\u00d723=\"),\n  @(\"17\u00d718=\", \"33\u00d750=\"),\n  @(\"91\u00d748=\", \"64\u00d776=\"),\n  @(\"32\u00d739=\", \"82\u00d720=\"),\n  @(\"78\u00d765=\", \"31\u00d715=\"),\n  @(\"95\u00d747=\", \"73\u00d745=\"),\n  @(\"36\u00d718=\", \"58\u00d770=\"),\n  @(\"26\u00d732=\", \"23\u00d739=\"),\n  @(\"78\u00d754=\", \"86\u00d748=\"),\n  @(\"58\u00d726=\", \"56\u00d750=\"),\n  @(\"15\u00d762=\", \"57\u00d711=\"),\n  @(\"79\u00d740=\", \"65\u00d763=\"),\n  @(\"53\u00d778=\", \"24\u00d797=\"),\n  @(\"98\u00d729=\", \"31\u00d738=\"),\n  @(\"20\u00d795=\", \"29\u00d770=\"),\n  @(\"74\u00d750=\", \"27\u00d797=\"),\n  @(\"44\u00d731=\", \"60\u00d788=\"),\n  @(\"26\u00d753=\", \"81\u00d763=\"),\n  @(\"54\u00d748=\", \"58\u00d767=\"),\n  @(\"59\u00d799=\", \"12\u00d789=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
